# "Add files via upload" — re-upload of Data/Tóm tắt dữ liệu.xlsx with the
# "window size 21" note bumped to "window size 31" (and the sheet view /
# column widths touched up as they were in the author's live Excel session).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Content: D1 held "window size 21" — update to "window size 31".
$ws.Range("D1").Value = "window size 31"

# Column widths as left by the author (B newly sized, C/D nudged, C no
# longer auto-bestFit since it was hand-set afterwards).
$ws.Columns.Item(2).ColumnWidth = 12.25
$ws.Columns.Item(3).ColumnWidth = 21.5
$ws.Columns.Item(4).ColumnWidth = 20.25

# View state: scrolled/selected a different cell before saving.
$excel.ActiveWindow.Zoom = 145
$ws.Range("D3").Select()
